$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new header cell for the standard deviation column
$ws.Range("I4").Value = "1 Yr Std Dev"

# Match the column width used by the author for the new column
# (real Excel would render 17+1/7 character-width units as 17.85546875 in the
# underlying OOXML "width" attribute; this is the closest this host's width
# model can reproduce)
$ws.Columns.Item(9).ColumnWidth = 17 + 1/7

# Match the active cell/selection recorded in the saved worksheet
$ws.Range("F8").Select() | Out-Null
